# "Cambios en crear plantillas" - add the new default "nucleo" row under the
# existing "Nombre Esquema" / "Propietario" header row of the schema-creation
# template sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "nucleo"
$ws.Range("B3").Value = "nucleo"

# Touch the (previously empty) remainder of the merged title cell so it picks
# up an explicit, plain format - matching how Excel materializes every cell of
# a merged range once the sheet is resaved.
$ws.Range("B1:C1").Style = "Normal"

# Leave the selection where the author last left it.
$null = $ws.Range("K13").Select()
